$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new trade record as row 16 (Principle, Start Principle, BuyPrice,
# SellPrice, IsShortSell, Price Change %, Date, Profitable)
$ws.Range("A16").Value = 9569.9699999999993
$ws.Range("B16").Value = 9622.9
$ws.Range("C16").Value = 78.48
$ws.Range("D16").Value = 78.05
$ws.Range("E16").Value = $false
$ws.Range("F16").Value = -0.55000000000000004
$ws.Range("G16").Value = 42624.611157407409
$ws.Range("H16").Value = $false
